$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2341597796143251
$ws.Range("C2").Value = 0.4958677685950413
$ws.Range("J2").Value = 0.02479338842975207
$ws.Range("P2").Value = 0.162534435261708
$ws.Range("S2").Value = 0.08264462809917356
$ws.Range("B3").Value = 0.01578947368421053
$ws.Range("C3").Value = 0.02631578947368421
$ws.Range("J3").Value = 0.04210526315789474
$ws.Range("P3").Value = 0.7263157894736842
$ws.Range("S3").Value = 0.1894736842105263
$ws.Range("J4").Value = 0.09433962264150944
$ws.Range("O4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.660377358490566
$ws.Range("S4").Value = 0.2264150943396226
$ws.Range("B6").Value = 0.06896551724137931
$ws.Range("D6").Value = 0.009852216748768473
$ws.Range("F6").Value = 0.05911330049261083
$ws.Range("J6").Value = 0.2315270935960591
$ws.Range("O6").Value = 0.01970443349753695
$ws.Range("Q6").Value = 0.1822660098522167
$ws.Range("R6").Value = 0.06403940886699508
$ws.Range("S6").Value = 0.3645320197044335
$ws.Range("B7").Value = 0.108433734939759
$ws.Range("D7").Value = 0.006024096385542169
$ws.Range("F7").Value = 0.07228915662650602
$ws.Range("J7").Value = 0.144578313253012
$ws.Range("O7").Value = 0.04216867469879518
$ws.Range("Q7").Value = 0.2108433734939759
$ws.Range("R7").Value = 0.07228915662650602
$ws.Range("S7").Value = 0.3433734939759036
$ws.Range("B8").Value = 0.09556313993174062
$ws.Range("D8").Value = 0.01706484641638225
$ws.Range("E8").Value = 0.003412969283276451
$ws.Range("F8").Value = 0.05802047781569966
$ws.Range("J8").Value = 0.1433447098976109
$ws.Range("O8").Value = 0.0204778156996587
$ws.Range("Q8").Value = 0.2320819112627986
$ws.Range("R8").Value = 0.1092150170648464
$ws.Range("S8").Value = 0.3208191126279863
$ws.Range("B9").Value = 0.1615720524017467
$ws.Range("D9").Value = 0.008733624454148471
$ws.Range("E9").Value = 0.004366812227074236
$ws.Range("F9").Value = 0.05240174672489083
$ws.Range("J9").Value = 0.1135371179039301
$ws.Range("O9").Value = 0.04366812227074236
$ws.Range("Q9").Value = 0.1703056768558952
$ws.Range("R9").Value = 0.06550218340611354
$ws.Range("S9").Value = 0.3799126637554585
$ws.Range("B10").Value = 0.1231732776617954
$ws.Range("D10").Value = 0.03061934585942937
$ws.Range("E10").Value = 0.00139178844815588
$ws.Range("F10").Value = 0.06541405706332637
$ws.Range("J10").Value = 0.1141266527487822
$ws.Range("O10").Value = 0.01809324982602644
$ws.Range("Q10").Value = 0.2386917188587335
$ws.Range("R10").Value = 0.08907446068197634
$ws.Range("S10").Value = 0.3194154488517745
$ws.Range("G11").Value = 0.1086142322097378
$ws.Range("J11").Value = 0.1310861423220974
$ws.Range("K11").Value = 0.1797752808988764
$ws.Range("L11").Value = 0.5730337078651685
$ws.Range("S11").Value = 0.00749063670411985
$ws.Range("G12").Value = 0.7329192546583851
$ws.Range("J12").Value = 0.1801242236024845
$ws.Range("K12").Value = 0.0124223602484472
$ws.Range("L12").Value = 0.04968944099378882
$ws.Range("S12").Value = 0.02484472049689441
$ws.Range("G13").Value = 0.71875
$ws.Range("J13").Value = 0.28125
$ws.Range("F15").Value = 0.02205882352941177
$ws.Range("H15").Value = 0.1102941176470588
$ws.Range("I15").Value = 0.06985294117647059
$ws.Range("J15").Value = 0.3639705882352941
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("M15").Value = 0.003676470588235294
$ws.Range("N15").Value = 0.003676470588235294
$ws.Range("O15").Value = 0.05514705882352941
$ws.Range("S15").Value = 0.3125
$ws.Range("F16").Value = 0.01345291479820628
$ws.Range("H16").Value = 0.1210762331838565
$ws.Range("I16").Value = 0.09865470852017937
$ws.Range("J16").Value = 0.4618834080717489
$ws.Range("K16").Value = 0.1165919282511211
$ws.Range("M16").Value = 0.008968609865470852
$ws.Range("O16").Value = 0.03587443946188341
$ws.Range("S16").Value = 0.1434977578475336
$ws.Range("F17").Value = 0.005825242718446602
$ws.Range("H17").Value = 0.09320388349514563
$ws.Range("I17").Value = 0.1320388349514563
$ws.Range("J17").Value = 0.516504854368932
$ws.Range("K17").Value = 0.06213592233009708
$ws.Range("M17").Value = 0.0116504854368932
$ws.Range("O17").Value = 0.07572815533980583
$ws.Range("S17").Value = 0.1029126213592233
$ws.Range("F18").Value = 0.01515151515151515
$ws.Range("H18").Value = 0.08080808080808081
$ws.Range("I18").Value = 0.1161616161616162
$ws.Range("J18").Value = 0.4141414141414141
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.0202020202020202
$ws.Range("N18").Value = 0.005050505050505051
$ws.Range("O18").Value = 0.1060606060606061
$ws.Range("S18").Value = 0.1313131313131313
$ws.Range("F19").Value = 0.01024765157984629
$ws.Range("H19").Value = 0.147736976942784
$ws.Range("I19").Value = 0.08198121263877028
$ws.Range("J19").Value = 0.426131511528608
$ws.Range("K19").Value = 0.09991460290350128
$ws.Range("M19").Value = 0.01707941929974381
$ws.Range("O19").Value = 0.09393680614859094
$ws.Range("S19").Value = 0.1229718189581554
